$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to be stored as text (not auto-converted to a number),
    # using an apostrophe prefix, then reset the cell style back to Normal
    # so no stray "quote prefix" style is left behind on the cell.
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextValue "D2" "41.675.25"
$ws.Range("E2").Value = "  +0.24%  "

Set-TextValue "D3" "2.478.77"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue "D5" "319.46"
$ws.Range("E5").Value = "  +1.63%  "

Set-TextValue "D6" "92.59"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("E8").Value = "  +0.07%  "

Set-TextValue "D9" "0.512"
$ws.Range("E9").Value = "  +0.72%  "

Set-TextValue "D10" "0.0866"
$ws.Range("E10").Value = "  +8.95%  "

Set-TextValue "D11" "33.13"
$ws.Range("E11").Value = "  +2.62%  "

$ws.Range("E12").Value = "  +0.14%  "

Set-TextValue "D13" "2.861.29"
$ws.Range("E13").Value = "  +0.93%  "

Set-TextValue "D14" "6.89"
$ws.Range("E14").Value = "  +1.03%  "

Set-TextValue "D15" "15.56"
$ws.Range("E15").Value = "  -1.17%  "

Set-TextValue "D16" "2.464.20"
$ws.Range("E16").Value = "  -0.86%  "

Set-TextValue "D17" "0.795"
$ws.Range("E17").Value = "  +2.66%  "

Set-TextValue "D18" "41.636.43"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("E19").Value = "  +0.03%  "

Set-TextValue "D20" "0.0₃0942"
$ws.Range("E20").Value = "  +0.97%  "

Set-TextValue "D21" "70.86"

$ws.Range("E22").Value = "  -0.64%  "

Set-TextValue "D23" "239.72"
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("E24").Value = "  +1.96%  "

$ws.Range("E25").Value = "  +2.74%  "

$ws.Range("E27").Value = "  +3.05%  "

$ws.Range("E28").Value = "  -0.54%  "

Set-TextValue "D29" "9.74"
$ws.Range("E29").Value = "  +0.84%  "

Set-TextValue "D30" "36.65"
$ws.Range("E30").Value = "  +4.85%  "

Set-TextValue "D31" "157.52"
$ws.Range("E31").Value = "  +1.29%  "

Set-TextValue "D32" "5.43"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("E33").Value = "  +0.04%  "

Set-TextValue "D34" "0.0765"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("E35").Value = "  -0.11%  "

Set-TextValue "D36" "17.17"
$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("E38").Value = "  +2.95%  "

$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("E40").Value = "  +1.98%  "

Set-TextValue "D41" "4.02"
$ws.Range("E41").Value = "  +2.05%  "

Set-TextValue "D42" "2.47"
$ws.Range("E42").Value = "  -0.63%  "

Set-TextValue "D43" "1.999.31"
$ws.Range("E43").Value = "  +1.61%  "

$ws.Range("E44").Value = "  +1.09%  "

Set-TextValue "D45" "18.66"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("E46").Value = "  +2.84%  "

Set-TextValue "D47" "9.57"
$ws.Range("E47").Value = "  +7.30%  "

Set-TextValue "D48" "2.719.33"
$ws.Range("E48").Value = "  +0.97%  "

Set-TextValue "D49" "98.05"
$ws.Range("E49").Value = "  +1.79%  "

Set-TextValue "D50" "75.66"
$ws.Range("E50").Value = "  +5.57%  "

Set-TextValue "D51" "67.16"
$ws.Range("E51").Value = "  +1.22%  "
